# Commit: "push - rename + update"
#
# A new pair of rows for client "Cocamar" (id Runrunit 157821) is inserted
# right after the header row, pushing all the pre-existing data rows down
# by two. As a knock-on effect of the new event being inserted earlier in
# the timeline, the Healthscore (column K) of what was the first
# pre-existing data row (Mart Minas / id 149896, "Metas não atingidas",
# now at row 4) is recalculated from 10 down to 8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 2 (the first data row),
# shifting all existing data down by two rows.
$ws.Rows("2:3").Insert()

# The inserted rows pick up row 1's (header) formatting by default;
# strip that back down to the plain/default style used by the other
# data rows.
$ws.Rows("2:3").ClearFormats()

# Column G ("data") holds a date-look-alike string ("yyyy-mm-dd"); format
# it as Text first so Excel's auto-detection doesn't silently convert the
# literal into a real date serial number, matching the other rows where
# the "data" column is plain text.
$ws.Range("G2:G3").NumberFormat = "@"

# Row 2: Cocamar / id 157821 - "Cliente pediu proposta"
$ws.Cells.Item(2, 1).Value = 157821
$ws.Cells.Item(2, 2).Value = "Cocamar, 15/07/2024"
$ws.Cells.Item(2, 3).Value = "backlog"
$ws.Cells.Item(2, 4).Value = "Acompanhamento de clientes"
$ws.Cells.Item(2, 5).Value = "[]"
$ws.Cells.Item(2, 6).Value = "Cliente pediu proposta"
$ws.Cells.Item(2, 7).Value = "2024-07-15"
$ws.Cells.Item(2, 8).Value = "Cocamar"
$ws.Cells.Item(2, 9).Value = "Cliente pediu proposta"
$ws.Cells.Item(2, 10).Value = 2.5
$ws.Cells.Item(2, 11).Value = 10

# Row 3: Cocamar / id 157821 - "Resolveu problema"
$ws.Cells.Item(3, 1).Value = 157821
$ws.Cells.Item(3, 2).Value = "Cocamar, 15/07/2024"
$ws.Cells.Item(3, 3).Value = "backlog"
$ws.Cells.Item(3, 4).Value = "Acompanhamento de clientes"
$ws.Cells.Item(3, 5).Value = "[]"
$ws.Cells.Item(3, 6).Value = "Resolveu problema"
$ws.Cells.Item(3, 7).Value = "2024-07-15"
$ws.Cells.Item(3, 8).Value = "Cocamar"
$ws.Cells.Item(3, 9).Value = "Resolveu problema"
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(3, 11).Value = 10

# The pre-existing first data row (Mart Minas / id 149896, "Metas não
# atingidas"), now shifted down to row 4, has its Healthscore recalculated
# from 10 to 8 because of the newly inserted preceding event.
$ws.Cells.Item(4, 11).Value = 8
